$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each changed cell as an explicit text value (matches source t="inlineStr" cells),
# avoiding Excel auto-converting numeric-looking strings into numbers.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "59.631.06"
Set-TextValue $ws.Range("E2") "  +3.06%  "
Set-TextValue $ws.Range("D3") "2.410.50"
Set-TextValue $ws.Range("E3") "  +3.28%  "
Set-TextValue $ws.Range("E4") "  +0.03%  "
Set-TextValue $ws.Range("D5") "552.51"
Set-TextValue $ws.Range("E5") "  +2.67%  "
Set-TextValue $ws.Range("D6") "137.19"
Set-TextValue $ws.Range("E6") "  +2.33%  "
Set-TextValue $ws.Range("E7") "  -0.02%  "
Set-TextValue $ws.Range("D8") "0.570"
Set-TextValue $ws.Range("E8") "  +2.34%  "
Set-TextValue $ws.Range("D9") "0.107"
Set-TextValue $ws.Range("E9") "  +6.16%  "
Set-TextValue $ws.Range("D10") "5.86"
Set-TextValue $ws.Range("E10") "  +6.09%  "
Set-TextValue $ws.Range("D11") "0.363"
Set-TextValue $ws.Range("E11") "  +2.05%  "
Set-TextValue $ws.Range("D12") "0.149"
Set-TextValue $ws.Range("E12") "  -2.48%  "
Set-TextValue $ws.Range("D13") "24.61"
Set-TextValue $ws.Range("E13") "  +4.04%  "
Set-TextValue $ws.Range("D14") "2.839.83"
Set-TextValue $ws.Range("D15") "59.523.13"
Set-TextValue $ws.Range("E15") "  +2.96%  "
Set-TextValue $ws.Range("D16") "0.0000140"
Set-TextValue $ws.Range("E16") "  +4.99%  "
Set-TextValue $ws.Range("D17") "2.447.61"
Set-TextValue $ws.Range("E17") "  +6.02%  "
Set-TextValue $ws.Range("D18") "11.32"
Set-TextValue $ws.Range("E18") "  +6.57%  "
Set-TextValue $ws.Range("D19") "4.41"
Set-TextValue $ws.Range("E19") "  +3.94%  "
Set-TextValue $ws.Range("D20") "336.50"
Set-TextValue $ws.Range("E20") "  +1.43%  "
Set-TextValue $ws.Range("E21") "  +5.08%  "
Set-TextValue $ws.Range("E22") "  +0.10%  "
Set-TextValue $ws.Range("D23") "64.79"
Set-TextValue $ws.Range("E23") "  +3.68%  "
Set-TextValue $ws.Range("E24") "  +0.77%  "
Set-TextValue $ws.Range("D25") "0.998"
Set-TextValue $ws.Range("E25") "  -0.12%  "
Set-TextValue $ws.Range("D26") "8.43"
Set-TextValue $ws.Range("E26") "  -0.82%  "
Set-TextValue $ws.Range("E27") "  -3.44%  "
Set-TextValue $ws.Range("E28") "  +2.51%  "
Set-TextValue $ws.Range("D29") "0.0₃0775"
Set-TextValue $ws.Range("E29") "  +5.88%  "
Set-TextValue $ws.Range("D30") "170.83"
Set-TextValue $ws.Range("E30") "  +0.06%  "
Set-TextValue $ws.Range("D31") "6.25"
Set-TextValue $ws.Range("E31") "  +2.72%  "
Set-TextValue $ws.Range("D32") "18.79"
Set-TextValue $ws.Range("E32") "  +1.80%  "
Set-TextValue $ws.Range("E33") "  +1.21%  "
Set-TextValue $ws.Range("D35") "4.32"
Set-TextValue $ws.Range("E35") "  +3.02%  "
Set-TextValue $ws.Range("E36") "  +4.91%  "
Set-TextValue $ws.Range("E37") "  +0.07%  "
Set-TextValue $ws.Range("E38") "  +1.30%  "
Set-TextValue $ws.Range("D39") "40.34"
Set-TextValue $ws.Range("E39") "  +3.59%  "
Set-TextValue $ws.Range("D40") "0.424"
Set-TextValue $ws.Range("E40") "  +13.23%  "
Set-TextValue $ws.Range("D41") "305.40"
Set-TextValue $ws.Range("E41") "  +7.16%  "
Set-TextValue $ws.Range("D42") "3.75"
Set-TextValue $ws.Range("E42") "  +3.40%  "
Set-TextValue $ws.Range("D43") "142.54"
Set-TextValue $ws.Range("E43") "  -0.98%  "
Set-TextValue $ws.Range("D44") "0.0961"
Set-TextValue $ws.Range("E44") "  +2.23%  "
Set-TextValue $ws.Range("D45") "0.0524"
Set-TextValue $ws.Range("E45") "  +4.45%  "
Set-TextValue $ws.Range("B46") "Polygon"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D46") "0.408"
Set-TextValue $ws.Range("E46") "  +5.96%  "
Set-TextValue $ws.Range("B47") "InjectiveProtocol"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D47") "19.15"
Set-TextValue $ws.Range("E47") "  +0.38%  "
Set-TextValue $ws.Range("B48") "Mantle"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D48") "0.571"
Set-TextValue $ws.Range("E48") "  +1.80%  "
Set-TextValue $ws.Range("B49") "VeChain"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D49") "0.0226"
Set-TextValue $ws.Range("E49") "  +4.15%  "
Set-TextValue $ws.Range("E50") "  -0.36%  "
Set-TextValue $ws.Range("E51") "  +4.67%  "
